{"js": "// The candidate's name in the cover letter header was updated again:\n// \"MYCO (MIRIAM) SULLIVAN\" -> \"MIRIAM SULLIVAN\".\n// Google Docs' own internal heading-anchor bookmark was also renamed\n// from \"_gjdgxs\" to \"_heading=h.gjdgxs\" when the document was re-saved.\n\nconst body = context.document.body;\n\n// 1) Rename the internal heading bookmark that sits at the very start\n//    of the title paragraph, preserving its (empty) location.\nconst OLD_BOOKMARK = \"_gjdgxs\";\nconst NEW_BOOKMARK = \"_heading=h.gjdgxs\";\n\nconst bookmarkRange = body.getBookmarkRangeOrNullObject(OLD_BOOKMARK);\nawait context.sync();\n\nif (!bookmarkRange.isNullObject) {\n  context.document.deleteBookmark(OLD_BOOKMARK);\n  bookmarkRange.insertBookmark(NEW_BOOKMARK);\n  await context.sync();\n}\n\n// 2) Replace the old name text with the new, shorter name while keeping\n//    the run's existing formatting (color/size/etc.) untouched.\nconst OLD_NAME = \"MYCO (MIRIAM) SULLIVAN\";\nconst NEW_NAME = \"MIRIAM SULLIVAN\";\n\nconst matches = body.search(OLD_NAME, { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < matches.items.length; i++) {\n  matches.items[i].insertText(NEW_NAME, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The candidate's name in the cover letter header was updated again:\n# \"MYCO (MIRIAM) SULLIVAN\" -> \"MIRIAM SULLIVAN\".\n# Google Docs' own internal heading-anchor bookmark was also renamed\n# from \"_gjdgxs\" to \"_heading=h.gjdgxs\" when the document was re-saved.\n\n$d = $word.ActiveDocument\n\n# 1) Rename the internal heading bookmark, preserving its (empty) location.\n$oldBookmarkName = \"_gjdgxs\"\n$newBookmarkName = \"_heading=h.gjdgxs\"\n\nif ($d.Bookmarks.Exists($oldBookmarkName)) {\n    $bm = $d.Bookmarks.Item($oldBookmarkName)\n    $bmStart = $bm.Range.Start\n    $bmEnd = $bm.Range.End\n    $bm.Delete()\n    $bmRange = $d.Range($bmStart, $bmEnd)\n    $d.Bookmarks.Add($newBookmarkName, $bmRange)\n}\n\n# 2) Replace the old name text with the new, shorter name while keeping\n#    the run's existing formatting (color/size/etc.) untouched.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"MYCO (MIRIAM) SULLIVAN\"\n$find.Replacement.Text = \"MIRIAM SULLIVAN\"\n$find.Execute([ref]$find.Text, $false, $true, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n"}
